$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master_table")
$ws.Rows("29:41").Delete()
